$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - new record appended to the Artfynd sheet
$ws.Range("A4").Value = 131242796
$ws.Range("B4").Value = 57881
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 100049
$ws.Range("F4").Value = "Spillkråka"
$ws.Range("G4").Value = "Dryocopus martius"
$ws.Range("H4").Value = "(Linnaeus, 1758)"
$ws.Range("M4").Value = "äldre spår"
$ws.Range("P4").Value = "Lille-Väktor, Boh"
$ws.Range("Q4").Value = 327429
$ws.Range("R4").Value = 6453589
$ws.Range("S4").Value = 5
$ws.Range("T4").Value = "Västra Götaland"
$ws.Range("U4").Value = "Lilla Edet"
$ws.Range("V4").Value = "Bohuslän"
$ws.Range("W4").Value = "Hjärtum"

# Date-looking values must stay as plain text, not get converted to
# Excel date serial numbers. Force text format while entering the value,
# then restore the default style so no stray number format lingers on
# the cell (matches the source file, where these are plain text cells).
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = "2026-01-31"
$ws.Range("Y4").Style = "Normal"

$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = "2026-01-31"
$ws.Range("AA4").Style = "Normal"

$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AG4").Value = $false

$ws.Range("AW4").Value = "Liv Vikingson"
$ws.Range("AX4").Value = "Liv Vikingson"

# These columns are blank for this record (no text entered), but the
# source row still carries a (typeless) cell placeholder at these
# addresses, so touch them without giving them any content.
$ws.Range("I4").Style = "Normal"
$ws.Range("K4").Style = "Normal"
$ws.Range("L4").Style = "Normal"
$ws.Range("N4").Style = "Normal"
$ws.Range("AT4").Style = "Normal"
$ws.Range("AY4").Style = "Normal"
